$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reassign the "Name" column (C) for Tables 1-4 (rows 1-20) to the new
# seating arrangement, keeping the Table/Seat labels in columns A/B intact.
$names = @(
    "Ivan ",
    "Alfiya",
    "Caroline",
    "Nasrin",
    "Gerrit",
    "Jens",
    "Viktor",
    "Alice",
    "Yanina",
    "Dario",
    "Niels",
    "Mahsa",
    "Ariana",
    "Nathalie",
    "Andrea",
    "Em",
    "Alexander",
    "Danil",
    "Afaf",
    "Miguel"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $names[$i]
}

# Remove the old Table 5 / Table 6 seating rows (old rows 21-30) entirely.
$ws.Range("A21:C30").ClearContents()

# Replace them with the "Remaining Colleagues" list (single column, A only).
$ws.Cells.Item(22, 1).Value = "Remaining Colleagues:"
$ws.Cells.Item(23, 1).Value = "Sweta"
$ws.Cells.Item(24, 1).Value = "Karel"
$ws.Cells.Item(25, 1).Value = "Fabienne"
$ws.Cells.Item(26, 1).Value = "Geraldine"
